# DFs-MRV-Limpo.xlsx edit
# Commit: "Melhor documentação e mais gráficos"
#
# The Balanco Patrimonial (BP) sheet previously had single line items
# (Contas a Receber, Estoques, Despesas Antecipadas, Outros, Financiamento,
# Debêntures, Outras Obrigações, Provisões) duplicated for both the
# "Circulante" (short term / CP) and "Não circulante" (long term / LP)
# sections. The edit disambiguates them with explicit "CP"/"LP" suffixes,
# and renames the total "Passivo" row to "Passivo Total".

$wb = $excel.ActiveWorkbook

$bp = $wb.Worksheets.Item("BP")

# --- Ativo Circulante (current assets) block ---
$bp.Range("A5").Value  = "Contas a Receber CP"
$bp.Range("A6").Value  = "Estoques CP"
$bp.Range("A8").Value  = "Despesas Antecipadas CP"
$bp.Range("A9").Value  = "Outros CP"

# --- Ativo Não Circulante (non-current assets) block ---
$bp.Range("A11").Value = "Contas a Receber LP"
$bp.Range("A12").Value = "Estoques LP"
$bp.Range("A13").Value = "Despesas Antecipadas LP"
$bp.Range("A18").Value = "Outros LP"

# --- Passivo Circulante (current liabilities) block ---
$bp.Range("A24").Value = "Financiamento CP"
$bp.Range("A25").Value = "Debêntures CP"
$bp.Range("A26").Value = "Provisões CP"
$bp.Range("A27").Value = "Outras Obrigações CP"

# --- Passivo Não Circulante (non-current liabilities) block ---
$bp.Range("A29").Value = "Financiamento LP"
$bp.Range("A30").Value = "Debêntures LP"
$bp.Range("A31").Value = "Outras Obrigações LP"
$bp.Range("A33").Value = "Provisões LP"

# --- Passivo total label ---
$bp.Range("A19").Value = "Passivo Total"

# Restore the cursor/selection on each sheet to match the saved view state.
$bp.Range("G17").Select()

$dfc = $wb.Worksheets.Item("DFC")
$dfc.Activate()
$dfc.Range("A8").Select()

$dre = $wb.Worksheets.Item("DRE")
$dre.Range("Q12").Select()

# DFC was the active tab when the workbook was last saved.
$dfc.Activate()

Write-Output "done"
